$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G16").Value = "15/10/2018"
$ws.Range("H16").Value = "30/10/2018"
$ws.Range("G18").Select()
